$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the test-scenario text in D8: the "special characters" test case
# is replaced with a "different languages" test case.
$ws.Range("D8").Value = "Verify the search field handles different languages"

# Move the active selection to E3, matching the saved cursor position.
$ws.Range("E3").Select()
